# Apply updated symbol list values scraped on 2022-12-29
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'244.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.204"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05748"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'3.237"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8149"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8695"
$ws.Range("D9").Style = "Normal"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1373"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.06935"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = "'0.03164"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.02942"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.09340"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = "'3.826"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = "'0.001531"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = "'0.04716"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = "'0.0005998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("D19").Value = "'0.006155"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'0.004102"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00008701"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.581"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'0.3186"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "'0.0002327"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03723"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006250"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1053"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002300"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = "'0.007517"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005435"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.4248"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.003309"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '47BOLOBOLOBestin24h'
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").Style = "Normal"
